$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Total pesquia" -> "Total pesquisa"
$ws.Range("E1").Value = "Total pesquisa"

# Update the active selection to E1, matching the saved selection state
$ws.Range("E1").Select()
